$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update user/admin values to lowercase
$ws.Range("A2").Value = "mauricio"
$ws.Range("C2").Value = "admin"

# Update the active selection to C2
$ws.Range("C2").Select()
